$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.910632418942356
$ws.Range("D2").Value = 0.0268843822672693
$ws.Range("E2").Value = 0.4467350302404096
$ws.Range("F2").Value = 0.5689882716264592
$ws.Range("G2").Value = 0.4091863476227573
$ws.Range("H2").Value = 0.5691357661397589
$ws.Range("K2").Value = 0.4151787958436444
$ws.Range("L2").Value = 0.09260658787896148
$ws.Range("M2").Value = 0.1755718693677686
$ws.Range("N2").Value = 1.912199208432453
$ws.Range("O2").Value = 1.90399042961792

$ws.Range("B3").Value = 0.8916484417135848
$ws.Range("D3").Value = 0.02430413097302875
$ws.Range("E3").Value = 0.4503219822598235
$ws.Range("F3").Value = 0.5668236749523956
$ws.Range("G3").Value = 0.4084925655182445
$ws.Range("H3").Value = 0.5719123362896497
$ws.Range("K3").Value = 0.3743562180654578
$ws.Range("L3").Value = 0.08692448544277198
$ws.Range("M3").Value = 0.1701456212863981
$ws.Range("N3").Value = 1.929860197451189
$ws.Range("O3").Value = 1.908105293697972

$ws.Range("B4").Value = 0.8803870732146493
$ws.Range("D4").Value = 0.02270806076525389
$ws.Range("E4").Value = 0.4526609001479356
$ws.Range("F4").Value = 0.5658280754547107
$ws.Range("G4").Value = 0.4083331038012972
$ws.Range("H4").Value = 0.5738472880450303
$ws.Range("K4").Value = 0.3491234353667494
$ws.Range("L4").Value = 0.08346945971012332
$ws.Range("M4").Value = 0.1668894609335965
$ws.Range("N4").Value = 1.941289529244072
$ws.Range("O4").Value = 1.911670501959506

$ws.Range("B5").Value = 0.8758978178163659
$ws.Range("D5").Value = 0.02205471701878992
$ws.Range("E5").Value = 0.4536484086437982
$ws.Range("F5").Value = 0.565506303258779
$ws.Range("G5").Value = 0.4083351506713555
$ws.Range("H5").Value = 0.5746937381793771
$ws.Range("K5").Value = 0.3387993720308629
$ws.Range("L5").Value = 0.08207010393743985
$ws.Range("M5").Value = 0.1655816795986524
$ws.Range("N5").Value = 1.946094368575695
$ws.Range("O5").Value = 1.913384704044148

$ws.Range("B6").Value = 0.8751584268722183
$ws.Range("D6").Value = 0.02194605361136581
$ws.Range("E6").Value = 0.4538144617225814
$ws.Range("F6").Value = 0.5654579457593982
$ws.Range("G6").Value = 0.408339539201485
$ws.Range("H6").Value = 0.5748377922980978
$ws.Range("K6").Value = 0.3370825813505292
$ws.Range("L6").Value = 0.08183826387134019
$ws.Range("M6").Value = 0.1653656826337482
$ws.Range("N6").Value = 1.946901107261329
$ws.Range("O6").Value = 1.913685135934074

$ws.Range("B7").Value = 0.8803261245522833
$ws.Range("D7").Value = 0.02269926137351774
$ws.Range("E7").Value = 0.4526740787568855
$ws.Range("F7").Value = 0.5658233959244114
$ws.Range("G7").Value = 0.4083328599967189
$ws.Range("H7").Value = 0.5738584688653532
$ws.Range("K7").Value = 0.3489843686290612
$ws.Range("L7").Value = 0.08345055258234169
$ws.Range("M7").Value = 0.166871746111827
$ws.Range("N7").Value = 1.941353732483634
$ws.Range("O7").Value = 1.911692561863887

$ws.Range("B8").Value = 0.9040051106618989
$ws.Range("D8").Value = 0.02599717585069072
$ws.Range("E8").Value = 0.4479435095427684
$ws.Range("F8").Value = 0.5681727576387487
$ws.Range("G8").Value = 0.4088918188784305
$ws.Range("H8").Value = 0.5700454090101843
$ws.Range("K8").Value = 0.4011383155416297
$ws.Range("L8").Value = 0.09064043659100207
$ws.Range("M8").Value = 0.173685286276946
$ws.Range("N8").Value = 1.918167268604781
$ws.Range("O8").Value = 1.905193757986481

$ws.Range("B9").Value = 0.9535516189847328
$ws.Range("D9").Value = 0.03236975910076723
$ws.Range("E9").Value = 0.4397476310559281
$ws.Range("F9").Value = 0.5754229252035401
$ws.Range("G9").Value = 0.4121029712428879
$ws.Range("H9").Value = 0.5643908952618517
$ws.Range("K9").Value = 0.5020592230901855
$ws.Range("L9").Value = 0.1050046544880274
$ws.Range("M9").Value = 0.1876415849453217
$ws.Range("N9").Value = 1.877339197256497
$ws.Range("O9").Value = 1.900685425126824

$ws.Range("B10").Value = 0.9918270235492344
$ws.Range("D10").Value = 0.03699293430884865
$ws.Range("E10").Value = 0.4343815364711716
$ws.Range("F10").Value = 0.5823585593567273
$ws.Range("G10").Value = 0.4157530239530871
$ws.Range("H10").Value = 0.5613439844915575
$ws.Range("K10").Value = 0.5753563133135913
$ws.Range("L10").Value = 0.1157160934327379
$ws.Range("M10").Value = 0.1982528572146407
$ws.Range("N10").Value = 1.850164803584363
$ws.Range("O10").Value = 1.902389325697129

$ws.Range("B11").Value = 1.009641587715464
$ws.Range("D11").Value = 0.0390831678556367
$ws.Range("E11").Value = 0.4320819461287115
$ws.Range("F11").Value = 0.5858627093527815
$ws.Range("G11").Value = 0.4176941809051868
$ws.Range("H11").Value = 0.5601975564679407
$ws.Range("K11").Value = 0.6085117166392422
$ws.Range("L11").Value = 0.1206226578250664
$ws.Range("M11").Value = 0.2031568261042693
$ws.Range("N11").Value = 1.838413753865765
$ws.Range("O11").Value = 1.904252895409741

$ws.Range("B12").Value = 1.016444935424033
$ws.Range("D12").Value = 0.03987280615241673
$ws.Range("E12").Value = 0.4312314409679292
$ws.Range("F12").Value = 0.587239771823306
$ws.Range("G12").Value = 0.41846962296691
$ws.Range("H12").Value = 0.55979782405268
$ws.Range("K12").Value = 0.6210392319623281
$ws.Range("L12").Value = 0.1224854357240872
$ws.Range("M12").Value = 0.2050247695506187
$ws.Range("N12").Value = 1.834051675339875
$ws.Range("O12").Value = 1.905114979525905

$ws.Range("B13").Value = 1.014977170900067
$ws.Range("D13").Value = 0.03970282796448998
$ws.Range("E13").Value = 0.4314137103286555
$ws.Range("F13").Value = 0.5869409691460845
$ws.Range("G13").Value = 0.4183008222842517
$ws.Range("H13").Value = 0.5598823848575449
$ws.Range("K13").Value = 0.6183424513811815
$ws.Range("L13").Value = 0.1220840427401555
$ws.Range("M13").Value = 0.2046219911103506
$ws.Range("N13").Value = 1.834987223037199
$ws.Range("O13").Value = 1.904922360830113

$ws.Range("B14").Value = 1.010200156971365
$ws.Range("D14").Value = 0.0391481699192866
$ws.Range("E14").Value = 0.4320115679692904
$ws.Range("F14").Value = 0.5859749971706236
$ws.Range("G14").Value = 0.4177571679610281
$ws.Range("H14").Value = 0.5601639812755224
$ws.Range("K14").Value = 0.6095429224057511
$ws.Range("L14").Value = 0.1207758147084519
$ws.Range("M14").Value = 0.2033102849129378
$ws.Range("N14").Value = 1.83805312379048
$ws.Range("O14").Value = 1.904320686027233

$ws.Range("B15").Value = 1.007281550711468
$ws.Range("D15").Value = 0.03880817891008803
$ws.Range("E15").Value = 0.4323804153729203
$ws.Range("F15").Value = 0.5853898356745333
$ws.Range("G15").Value = 0.4174294210690164
$ws.Range("H15").Value = 0.5603409445855334
$ws.Range("K15").Value = 0.6041493262816573
$ws.Range("L15").Value = 0.1199751051945839
$ws.Range("M15").Value = 0.202508244714771
$ws.Range("N15").Value = 1.839942508348322
$ws.Range("O15").Value = 1.903972505687335

$ws.Range("B16").Value = 0.9906708578273538
$ws.Range("D16").Value = 0.03685607047158612
$ws.Range("E16").Value = 0.4345346631517764
$ws.Range("F16").Value = 0.5821365725735248
$ws.Range("G16").Value = 0.4156318137583668
$ws.Range("H16").Value = 0.5614237227734691
$ws.Range("K16").Value = 0.5731856898641468
$ws.Range("L16").Value = 0.1153961108037009
$ws.Range("M16").Value = 0.1979339062609213
$ws.Range("N16").Value = 1.850945048777035
$ws.Range("O16").Value = 1.902289428669604

$ws.Range("B17").Value = 0.9805835158946934
$ws.Range("D17").Value = 0.03565519221504587
$ws.Range("E17").Value = 0.4358924277266274
$ws.Range("F17").Value = 0.580230154351618
$ws.Range("G17").Value = 0.4146009414517096
$ws.Range("H17").Value = 0.5621493014384527
$ws.Range("K17").Value = 0.5541418943303427
$ws.Range("L17").Value = 0.1125956495791485
$ws.Range("M17").Value = 0.1951472852674101
$ws.Range("N17").Value = 1.857851147308914
$ws.Range("O17").Value = 1.901535604888892

$ws.Range("B18").Value = 0.9748195000838962
$ws.Range("D18").Value = 0.03496326943858463
$ws.Range("E18").Value = 0.4366866967124556
$ws.Range("F18").Value = 0.579166496521772
$ws.Range("G18").Value = 0.4140344346684941
$ws.Range("H18").Value = 0.5625891936855538
$ws.Range("K18").Value = 0.5431707716185485
$ws.Range("L18").Value = 0.1109880932033036
$ws.Range("M18").Value = 0.1935517358772678
$ws.Range("N18").Value = 1.861880837947545
$ws.Range("O18").Value = 1.90120448387907

$ws.Range("B19").Value = 0.9728744399797904
$ws.Range("D19").Value = 0.03472878968359083
$ws.Range("E19").Value = 0.4369579113858126
$ws.Range("F19").Value = 0.5788120070914857
$ws.Range("G19").Value = 0.4138471635015293
$ws.Range("H19").Value = 0.5627420100246923
$ws.Range("K19").Value = 0.5394531307268551
$ws.Range("L19").Value = 0.1104443540892532
$ws.Range("M19").Value = 0.1930127580418883
$ws.Range("N19").Value = 1.863255096138101
$ws.Range("O19").Value = 1.901109973208321

$ws.Range("B20").Value = 0.9816534064839573
$ws.Range("D20").Value = 0.03578315318591052
$ws.Range("E20").Value = 0.4357465132518472
$ws.Range("F20").Value = 0.5804296951703165
$ws.Range("G20").Value = 0.4147079448151771
$ws.Range("H20").Value = 0.5620697280784128
$ws.Range("K20").Value = 0.5561709691911574
$ws.Range("L20").Value = 0.1128934335479244
$ws.Range("M20").Value = 0.1954431773501746
$ws.Range("N20").Value = 1.85711003204101
$ws.Range("O20").Value = 1.901605247568483

$ws.Range("B21").Value = 1.011601729651375
$ws.Range("D21").Value = 0.03931113796004126
$ws.Range("E21").Value = 0.4318354121054568
$ws.Range("F21").Value = 0.5862573670606466
$ws.Range("G21").Value = 0.4179157568460568
$ws.Range("H21").Value = 0.5600803366137228
$ws.Range("K21").Value = 0.6121283152380386
$ws.Range("L21").Value = 0.1211599444015832
$ws.Range("M21").Value = 0.2036952696063139
$ws.Range("N21").Value = 1.83715021209272
$ws.Range("O21").Value = 1.904493169133787

$ws.Range("B22").Value = 1.031508753637098
$ws.Range("D22").Value = 0.04160585928539717
$ws.Range("E22").Value = 0.4293975839979041
$ws.Range("F22").Value = 0.5903581680800869
$ws.Range("G22").Value = 0.4202475423555825
$ws.Range("H22").Value = 0.5589806054262283
$ws.Range("K22").Value = 0.6485377941978925
$ws.Range("L22").Value = 0.1265903332599692
$ws.Range("M22").Value = 0.2091520584767181
$ws.Range("N22").Value = 1.824617004306983
$ws.Range("O22").Value = 1.907292096625469

$ws.Range("B23").Value = 1.020853624046708
$ws.Range("D23").Value = 0.04038214393172268
$ws.Range("E23").Value = 0.4306878884597669
$ws.Range("F23").Value = 0.5881427956071761
$ws.Range("G23").Value = 0.418981495493199
$ws.Range("H23").Value = 0.5595492321877771
$ws.Range("K23").Value = 0.6291204166690534
$ws.Range("L23").Value = 0.1236895278748023
$ws.Range("M23").Value = 0.2062338952080083
$ws.Range("N23").Value = 1.831259403452762
$ws.Range("O23").Value = 1.905714896551501

$ws.Range("B24").Value = 0.9811695988517215
$ws.Range("D24").Value = 0.03572530676881058
$ws.Range("E24").Value = 0.4358124385962485
$ws.Range("F24").Value = 0.5803393819387637
$ws.Range("G24").Value = 0.414659487122691
$ws.Range("H24").Value = 0.5621056323333136
$ws.Range("K24").Value = 0.5552536947616886
$ws.Range("L24").Value = 0.11275879772667
$ws.Range("M24").Value = 0.1953093842326687
$ws.Range("N24").Value = 1.857444905590381
$ws.Range("O24").Value = 1.901573443565127

$ws.Range("B25").Value = 0.9398172458850524
$ws.Range("D25").Value = 0.03065604685910728
$ws.Range("E25").Value = 0.4418494878447037
$ws.Range("F25").Value = 0.573178842696187
$ws.Range("G25").Value = 0.4110076885581222
$ws.Range("H25").Value = 0.565725833986626
$ws.Range("K25").Value = 0.4749048340966056
$ws.Range("L25").Value = 0.1010907494194981
$ws.Range("M25").Value = 0.1838028704449215
$ws.Range("N25").Value = 1.887888206151567
$ws.Range("O25").Value = 1.901023818182878

Write-Output "applied changes"
